$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (A6) with the new "Test Onotlogy" entry
$rng = $ws.Range("A6")
$rng.Value = "Test Onotlogy"

# Match the look of the other "Ontology IRI" cells in column A (blue hyperlink-style
# font) combined with a boxed medium left/right border (no top/bottom) and
# center-wrapped text, as used elsewhere in this sheet.
$rng.Font.Color = 9516568
$rng.Borders.Item(7).Weight = -4138
$rng.Borders.Item(10).Weight = -4138
$rng.VerticalAlignment = -4108
$rng.WrapText = $true

# Move the active selection, matching the saved view state of the edited workbook
$null = $ws.Range("A13").Select()
